$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 51, shifting the existing rows 51-52 down to 52-53
$ws.Rows.Item(51).Insert()

# Populate the newly inserted row 51 with the new weekly record
$ws.Range("A51").Value = 11
$ws.Range("B51").Value = "Vega Monumental Concepción"
$ws.Range("C51").Value = "Bíobío"
$ws.Range("D51").Value = 45077
$ws.Range("D51").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E51").Value = 8
$ws.Range("F51").Value = 100112026
$ws.Range("G51").Value = "Haba"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 100
$ws.Range("K51").Value = 18000
$ws.Range("L51").Value = 20000
$ws.Range("M51").Value = 19000
$ws.Range("N51").Value = "$/saco 25 kilos"
$ws.Range("O51").Value = "Provincia de Limarí"
$ws.Range("P51").Value = 760
$ws.Range("Q51").Value = 25
$ws.Range("R51").Value = "Hortaliza"
